$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying naive-forecaster data rolled forward by one period: the
# oldest observation (old row 2) is dropped, every remaining row shifts up
# by one, and the lag/forecast columns (C and E) are recomputed with the
# bugfixed values for the resulting data vector.
$ws.Rows(2).Delete()

    $ws.Range("E2").Value = 1.782259294303912
    $ws.Range("E3").Value = -0.08289353495386509
    $ws.Range("C4").Value = -0.9921462019007898
    $ws.Range("E4").Value = -0.3230348957779294
    $ws.Range("C5").Value = -1.120856461469888
    $ws.Range("E5").Value = -0.363786394693788
    $ws.Range("C6").Value = 0.2467309912830284
    $ws.Range("E6").Value = -0.06959526544320083
    $ws.Range("C7").Value = 0.2946875655135184
    $ws.Range("E7").Value = -0.2139598932957232
    $ws.Range("C8").Value = 1.160201558804674
    $ws.Range("E8").Value = 0.2932139896134167
    $ws.Range("C9").Value = 1.242549344471122
    $ws.Range("E9").Value = 0.903223459378788
    $ws.Range("C10").Value = 1.048604932640185
    $ws.Range("E10").Value = 1.078804187516891
    $ws.Range("C11").Value = 1.180122656701199
    $ws.Range("E11").Value = 1.31837503023402
    $ws.Range("C12").Value = 1.06837811337479
    $ws.Range("E12").Value = 1.25598608434605
    $ws.Range("C13").Value = 1.381744454873757
    $ws.Range("E13").Value = 1.658305347589661
    $ws.Range("C14").Value = 1.435208340819005
    $ws.Range("E14").Value = 1.407107513712802
    $ws.Range("C15").Value = 1.404039427736437
    $ws.Range("E15").Value = 1.488472133572305
    $ws.Range("C16").Value = 1.577589817310243
    $ws.Range("E16").Value = 1.464859320654099
    $ws.Range("C17").Value = 1.593617458167307
    $ws.Range("E17").Value = 1.644157643645183
    $ws.Range("C18").Value = 1.979569114089963
    $ws.Range("E18").Value = 1.639881111696151
    $ws.Range("C19").Value = 1.885212754467758
    $ws.Range("E19").Value = 1.586470485311331
    $ws.Range("C20").Value = 1.815212363528707
    $ws.Range("E20").Value = 1.806931013599544
    $ws.Range("C21").Value = 1.877372574190161
    $ws.Range("E21").Value = 1.974604558490256
    $ws.Range("C22").Value = 2.140635848901895
    $ws.Range("E22").Value = 1.99288634244883
    $ws.Range("C23").Value = 2.284026378382942
    $ws.Range("E23").Value = 2.119133965447961
    $ws.Range("C24").Value = 2.336516087993035
    $ws.Range("E24").Value = 2.162438527487853
    $ws.Range("C25").Value = 2.279995067217899
    $ws.Range("E25").Value = 1.93172124148786
    $ws.Range("C26").Value = 1.404530461900833
    $ws.Range("E26").Value = 1.887821778955101
    $ws.Range("C27").Value = 1.264761787657309
    $ws.Range("E27").Value = 1.798687504247187
    $ws.Range("C28").Value = 1.17909021197069
    $ws.Range("E28").Value = 1.636329093826605
    $ws.Range("C29").Value = 1.025257057800411
    $ws.Range("E29").Value = 1.008270799755984
    $ws.Range("C30").Value = 0.5345697479163913
    $ws.Range("E30").Value = 1.247274949485733
    $ws.Range("C31").Value = 0.824608016336259
    $ws.Range("E31").Value = 1.395219579261608
    $ws.Range("C32").Value = -1.788000783651811
    $ws.Range("E32").Value = -0.02261741485058977
    $ws.Range("C33").Value = -1.788000783651811
    $ws.Range("E33").Value = -1.119700950349478
    $ws.Range("C34").Value = -2.680286313062752
    $ws.Range("E34").Value = -2.013357217277445
    $ws.Range("C35").Value = -1.4191429117966
    $ws.Range("E35").Value = -0.1125839228000469
    $ws.Range("C36").Value = -1.098964423305859
    $ws.Range("E36").Value = 1.055324027461602
    $ws.Range("C37").Value = -1.098964423305859
    $ws.Range("E37").Value = 0.5759895884974942
    $ws.Range("C38").Value = 1.514644056931957
    $ws.Range("E38").Value = -0.06175132635745095
    $ws.Range("C39").Value = 1.896944139870205
    $ws.Range("E39").Value = 0.5116467003986713
    $ws.Range("C40").Value = 1.916393754370604
    $ws.Range("E40").Value = 0.4136280550221194
    $ws.Range("C41").Value = 1.916393754370604
    $ws.Range("E41").Value = 0.3530477102890783
    $ws.Range("C42").Value = -0.8557279162653919
    $ws.Range("E42").Value = -0.2004689067778398
    $ws.Range("C43").Value = -0.7016063587211741
    $ws.Range("E43").Value = 0.2942159770784825
    $ws.Range("C44").Value = -0.7359525160776204
    $ws.Range("E44").Value = 0.6923809915882817
    $ws.Range("C45").Value = -0.7359525160776204
    $ws.Range("E45").Value = -0.01286797263981843
    $ws.Range("C46").Value = -0.1316183744203947
    $ws.Range("E46").Value = -0.09133135081734745
    $ws.Range("C47").Value = -0.1754728623905355
    $ws.Range("E47").Value = -0.0331361487157622
    $ws.Range("C48").Value = -0.187152549496028
    $ws.Range("E48").Value = -0.3126391654689975
    $ws.Range("C49").Value = -0.187152549496028
    $ws.Range("E49").Value = -0.1152140120150968
    $ws.Range("C50").Value = 0.3903331526556864
    $ws.Range("E50").Value = -0.2264357368625403
    $ws.Range("C51").Value = 0.5695821893874298
    $ws.Range("E51").Value = 0.316149716722669
    $ws.Range("C52").Value = 0.6150340712028246
    $ws.Range("E52").Value = 0.6473947787101642
